$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("is_active") currently stores a volatile =TRUE() formula (numeric 1).
# Replace it with the literal text value "TRUE" for rows 2-9, matching the
# fix described in the commit ("Fix: boolean values").
for ($r = 2; $r -le 9; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    # Build the text "TRUE" via a formula that returns a string, then freeze
    # it into a static value so the stored type is text (t="s") rather than
    # a boolean (t="b") or a live formula (t="str").
    $cell.Formula = "=""TRUE"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
$excel.CutCopyMode = $false

# Restore the selection/scroll position recorded in the saved view.
$ws.Range("G10").Select()
$excel.ActiveWindow.DisplayGridlines = $true
